$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "usual"
$ws.Range("B1").Value = "slow_release"
$ws.Range("B1").Select()
